# Update column G ("K") values on Sheet1 per regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals"). Only the K column values change; all other
# cells/columns are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 2
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 2
$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 2
$ws.Range("G25").Value = 2
$ws.Range("G26").Value = 5
$ws.Range("G27").Value = 2
$ws.Range("G28").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("G32").Value = 2
$ws.Range("G33").Value = 1
$ws.Range("G34").Value = 1
$ws.Range("G35").Value = 0
$ws.Range("G36").Value = 3
$ws.Range("G37").Value = 1
$ws.Range("G38").Value = 2
$ws.Range("G39").Value = 0
$ws.Range("G40").Value = 1
$ws.Range("G41").Value = 3
$ws.Range("G42").Value = 1
$ws.Range("G43").Value = 1
$ws.Range("G44").Value = 3
$ws.Range("G45").Value = 1
$ws.Range("G46").Value = 2
$ws.Range("G47").Value = 1
$ws.Range("G48").Value = 2
$ws.Range("G49").Value = 1
$ws.Range("G50").Value = 2
$ws.Range("G51").Value = 0
$ws.Range("G52").Value = 3
$ws.Range("G53").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("G56").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("G58").Value = 2
$ws.Range("G59").Value = 1
$ws.Range("G60").Value = 0
$ws.Range("G61").Value = 1
$ws.Range("G62").Value = 1
$ws.Range("G63").Value = 2
$ws.Range("G64").Value = 1
$ws.Range("G65").Value = 0
$ws.Range("G66").Value = 1
$ws.Range("G67").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("G70").Value = 1
$ws.Range("G71").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("G73").Value = 2
$ws.Range("G75").Value = 2
$ws.Range("G76").Value = 1
$ws.Range("G77").Value = 3
$ws.Range("G78").Value = 0
$ws.Range("G79").Value = 2
